$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the now-unused duplicate "Identifier"/"Identifier Type" columns (S:T) in rows 2-3
$ws.Range("S2:T3").Clear()

# R3 becomes "Catalog Key" (was "BIB")
$ws.Range("R3").Value = "Catalog Key"

# Select R3 to match the new selection state
$ws.Range("R3").Select()
